# Add an "EmployeeID" column (new column A) to the employee sheet, give the
# first two existing employees (Suhrob1 / renamed-to-admin Suhrob2) IDs, and
# blank out the remaining placeholder sample rows (4-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift existing header/data one column to the right is NOT needed -- the
# sheet already uses columns B..J; we just need to populate the now-used
# column A.

# --- Header row ---------------------------------------------------------
$ws.Range("A1").Value = "EmployeeID"
# match the look of the rest of the header row (B1 already carries the
# header style)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# --- Row 2: existing employee, now with an EmployeeID, moved to the admin
#     department / marketing role --------------------------------------
$ws.Range("A2").Value = "DK1003"
$ws.Range("G2").Value = "経営管理部"
$ws.Range("H2").Value = "Marketolog"

# --- Row 3: second employee becomes the admin, keeps default dept/role -
$ws.Range("A3").Value = "DK1004"
$ws.Range("B3").Value = "Suhrob2"
$ws.Range("G3").Value = "人事部"
$ws.Range("H3").Value = "CEO"

# --- Rows 4-11: clear out the remaining placeholder sample rows --------
for ($r = 4; $r -le 11; $r++) {
    $ws.Range("B${r}:J${r}").ClearContents()
}

# --- Selection matches the authored workbook ----------------------------
$ws.Range("A3").Select()
